$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Fix typo in monthly cost figures: -1028.97 -> -1848.19, rows 3-32, column C
for ($r = 3; $r -le 32; $r++) {
    $ws1.Cells.Item($r, 3).Value = -1848.19
}

# Make Sheet1 the active / selected sheet, set its view position + selection
$ws1.Activate()
$ws1.Range("D10").Select()

$excel.ActiveWindow.ScrollColumn = 3
